$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 13822.947
$ws.Range("I40").Value = 16214.286
$ws.Range("J40").Value = 12428
$ws.Range("K40").Value = 16214.286
$ws.Range("L40").Value = 12428
$ws.Range("M40").Value = -16039.286
$ws.Range("N40").Value = -12778

$ws.Range("H64").Value = 3311
$ws.Range("I64").Value = 3333.3333
$ws.Range("J64").Value = 3288.6667
$ws.Range("K64").Value = 3333.3333
$ws.Range("L64").Value = 3288.6667
$ws.Range("M64").Value = -3085.3333
$ws.Range("N64").Value = -3784.6667

$ws.Range("H67").Value = 3311
$ws.Range("I67").Value = 3333.3333
$ws.Range("J67").Value = 3288.6667
$ws.Range("K67").Value = 3333.3333
$ws.Range("L67").Value = 3288.6667
$ws.Range("M67").Value = -2475.3333
$ws.Range("N67").Value = -5004.6667

$ws.Range("H98").Value = 1204.96
$ws.Range("I98").Value = 1276.381
$ws.Range("K98").Value = 1276.381
$ws.Range("M98").Value = 221.6189999999999

$ws.Range("H122").Value = 1204.96
$ws.Range("I122").Value = 1276.381
$ws.Range("K122").Value = 3829.143
$ws.Range("M122").Value = -1379.143

$ws.Range("H129").Value = 1929.1177
$ws.Range("J129").Value = 2400.5833
$ws.Range("L129").Value = 7201.749899999999
$ws.Range("N129").Value = -17201.7499

$ws.Range("H137").Value = 200818.22
$ws.Range("I137").Value = 307022.16
$ws.Range("J137").Value = 6111
$ws.Range("K137").Value = 921066.48
$ws.Range("L137").Value = 18333
$ws.Range("M137").Value = -918516.48
$ws.Range("N137").Value = -23433

$ws.Range("H138").Value = 6548.1875
$ws.Range("I138").Value = 2462
$ws.Range("J138").Value = 8999.9
$ws.Range("K138").Value = 7386
$ws.Range("L138").Value = 26999.7
$ws.Range("M138").Value = -2246
$ws.Range("N138").Value = -37279.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19025.797
$ws.Range("I32").Value = 18139.297
$ws.Range("K32").Value = 18139.297
$ws.Range("M32").Value = -17852.297

$ws.Range("H45").Value = 56438.953
$ws.Range("I45").Value = 77243.39999999999
$ws.Range("K45").Value = 77243.39999999999
$ws.Range("M45").Value = -76866.39999999999

$ws.Range("H46").Value = 69057.60000000001
$ws.Range("I46").Value = 90000
$ws.Range("J46").Value = 63822
$ws.Range("K46").Value = 90000
$ws.Range("L46").Value = 63822
$ws.Range("M46").Value = -89681
$ws.Range("N46").Value = -64460

$ws.Range("H61").Value = 5419.2573
$ws.Range("I61").Value = 5850.05
$ws.Range("K61").Value = 5850.05
$ws.Range("M61").Value = -5638.05

$ws.Range("H74").Value = 1455.6
$ws.Range("I74").Value = 1360
$ws.Range("J74").Value = 1599
$ws.Range("K74").Value = 1360
$ws.Range("L74").Value = 1599
$ws.Range("M74").Value = -486
$ws.Range("N74").Value = -3347

$ws.Range("H77").Value = 1455.6
$ws.Range("I77").Value = 1360
$ws.Range("J77").Value = 1599
$ws.Range("K77").Value = 6800
$ws.Range("L77").Value = 7995
$ws.Range("M77").Value = -2432
$ws.Range("N77").Value = -16731

$ws.Range("H132").Value = 2984.8472
$ws.Range("I132").Value = 1236.0962
$ws.Range("J132").Value = 7531.6
$ws.Range("K132").Value = 3708.2886
$ws.Range("L132").Value = 22594.8
$ws.Range("M132").Value = -1178.2886
$ws.Range("N132").Value = -27654.8

$ws.Range("H136").Value = 5419.2573
$ws.Range("I136").Value = 5850.05
$ws.Range("K136").Value = 17550.15
$ws.Range("M136").Value = -15000.15

$ws.Range("H138").Value = 85000
$ws.Range("J138").Value = 85000
$ws.Range("L138").Value = 85000
$ws.Range("N138").Value = -95280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2833.2856
$ws.Range("I105").Value = 2899.5
$ws.Range("K105").Value = 2899.5
$ws.Range("M105").Value = -1152.5

$ws.Range("H134").Value = 5028.4727
$ws.Range("I134").Value = 2164.6453
$ws.Range("K134").Value = 6493.9359
$ws.Range("M134").Value = -3958.9359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4176.7
$ws.Range("I99").Value = 2329.6155
$ws.Range("K99").Value = 2329.6155
$ws.Range("M99").Value = -831.6154999999999

$ws.Range("H122").Value = 6133.353
$ws.Range("I122").Value = 4326.7
$ws.Range("J122").Value = 8714.286
$ws.Range("K122").Value = 12980.1
$ws.Range("L122").Value = 26142.858
$ws.Range("M122").Value = -10530.1
$ws.Range("N122").Value = -31042.858

$ws.Range("H126").Value = 4176.7
$ws.Range("I126").Value = 2329.6155
$ws.Range("K126").Value = 6988.8465
$ws.Range("M126").Value = -4518.8465

$ws.Range("H132").Value = 57972748
$ws.Range("I132").Value = 63493776
$ws.Range("K132").Value = 190481328
$ws.Range("M132").Value = -190478798

$ws.Range("H134").Value = 1759.44
$ws.Range("I134").Value = 1603.5834
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 4810.7502
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -2275.7502
$ws.Range("N134").Value = -21570

$ws.Range("H135").Value = 64843.75
$ws.Range("J135").Value = 64843.75
$ws.Range("L135").Value = 64843.75
$ws.Range("N135").Value = -74983.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76.21429000000001
$ws.Range("I2").Value = 49
$ws.Range("J2").Value = 91.333336
$ws.Range("K2").Value = 294
$ws.Range("L2").Value = 548.000016
$ws.Range("M2").Value = -181
$ws.Range("N2").Value = -774.000016

$ws.Range("H14").Value = 262.13043
$ws.Range("I14").Value = 262.13043
$ws.Range("K14").Value = 786.39129
$ws.Range("M14").Value = -613.39129

$ws.Range("H38").Value = 325.33334
$ws.Range("I38").Value = 109.4
$ws.Range("J38").Value = 521.63635
$ws.Range("K38").Value = 328.2
$ws.Range("L38").Value = 1564.90905
$ws.Range("M38").Value = 18.79999999999995
$ws.Range("N38").Value = -2258.90905

$ws.Range("H50").Value = 263.33334
$ws.Range("I50").Value = 232.83333
$ws.Range("J50").Value = 324.33334
$ws.Range("K50").Value = 698.49999
$ws.Range("L50").Value = 973.0000200000001
$ws.Range("M50").Value = -217.49999
$ws.Range("N50").Value = -1935.00002

$ws.Range("H53").Value = 263.33334
$ws.Range("I53").Value = 232.83333
$ws.Range("J53").Value = 324.33334
$ws.Range("K53").Value = 698.49999
$ws.Range("L53").Value = 973.0000200000001
$ws.Range("M53").Value = -217.49999
$ws.Range("N53").Value = -1935.00002

$ws.Range("H128").Value = 145399.2
$ws.Range("I128").Value = 145399.2
$ws.Range("K128").Value = 436197.6
$ws.Range("M128").Value = -431217.6

$ws.Range("H138").Value = 3029.68
$ws.Range("I138").Value = 3064.1904
$ws.Range("J138").Value = 2848.5
$ws.Range("K138").Value = 9192.5712
$ws.Range("L138").Value = 8545.5
$ws.Range("M138").Value = -4052.5712
$ws.Range("N138").Value = -18825.5

$ws.Range("H139").Value = 2391.6428
$ws.Range("I139").Value = 1495
$ws.Range("K139").Value = 4485
$ws.Range("M139").Value = 655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 276.6111
$ws.Range("J2").Value = 357.2857
$ws.Range("L2").Value = 357.2857
$ws.Range("N2").Value = -583.2857

$ws.Range("H70").Value = 5094.75
$ws.Range("I70").Value = 5053
$ws.Range("J70").Value = 5164.3335
$ws.Range("K70").Value = 5053
$ws.Range("L70").Value = 5164.3335
$ws.Range("M70").Value = -4783
$ws.Range("N70").Value = -5704.3335

$ws.Range("H73").Value = 5094.75
$ws.Range("I73").Value = 5053
$ws.Range("J73").Value = 5164.3335
$ws.Range("K73").Value = 5053
$ws.Range("L73").Value = 5164.3335
$ws.Range("M73").Value = -4117
$ws.Range("N73").Value = -7036.3335

$ws.Range("H97").Value = 1489.8334
$ws.Range("I97").Value = 1054.7333
$ws.Range("K97").Value = 1054.7333
$ws.Range("M97").Value = -558.7333000000001

$ws.Range("H122").Value = 362809.53
$ws.Range("J122").Value = 6940.8125
$ws.Range("L122").Value = 20822.4375
$ws.Range("N122").Value = -25722.4375

$ws.Range("H132").Value = 56335.566
$ws.Range("I132").Value = 82010.32000000001
$ws.Range("J132").Value = 2846.5
$ws.Range("K132").Value = 246030.96
$ws.Range("L132").Value = 8539.5
$ws.Range("M132").Value = -243500.96
$ws.Range("N132").Value = -13599.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5060.1
$ws.Range("I7").Value = 3636.3333
$ws.Range("J7").Value = 5670.2856
$ws.Range("K7").Value = 3636.3333
$ws.Range("L7").Value = 5670.2856
$ws.Range("M7").Value = -3524.3333
$ws.Range("N7").Value = -5894.2856

$ws.Range("H16").Value = 1147.8
$ws.Range("I16").Value = 632.75
$ws.Range("J16").Value = 3208
$ws.Range("K16").Value = 632.75
$ws.Range("L16").Value = 3208
$ws.Range("M16").Value = -462.75
$ws.Range("N16").Value = -3548

$ws.Range("H126").Value = 5060.1
$ws.Range("I126").Value = 3636.3333
$ws.Range("J126").Value = 5670.2856
$ws.Range("K126").Value = 10908.9999
$ws.Range("L126").Value = 17010.8568
$ws.Range("M126").Value = -8438.999899999999
$ws.Range("N126").Value = -21950.8568

$ws.Range("H132").Value = 2874.5
$ws.Range("J132").Value = 2909.4482
$ws.Range("L132").Value = 8728.3446
$ws.Range("N132").Value = -13788.3446

$ws.Range("H136").Value = 5893
$ws.Range("I136").Value = 3013.6875
$ws.Range("K136").Value = 9041.0625
$ws.Range("M136").Value = -6491.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 392.3846
$ws.Range("I4").Value = 341.75
$ws.Range("K4").Value = 341.75
$ws.Range("M4").Value = -228.75

$ws.Range("H126").Value = 500002500
$ws.Range("I126").Value = 1000000000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 3000000000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -2999997530
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 2206.3333
$ws.Range("I132").Value = 1137
$ws.Range("K132").Value = 3411
$ws.Range("M132").Value = -881
